$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 152.66667
$ws.Range("I4").Value = 152.66667
$ws.Range("K4").Value = 152.66667
$ws.Range("M4").Value = -38.66667000000001

$ws.Range("H8").Value = 268.9
$ws.Range("I8").Value = 115.625
$ws.Range("K8").Value = 346.875
$ws.Range("M8").Value = -207.875

$ws.Range("H9").Value = 1467.3334
$ws.Range("I9").Value = 1280.9
$ws.Range("J9").Value = 2399.5
$ws.Range("K9").Value = 1280.9
$ws.Range("L9").Value = 2399.5
$ws.Range("M9").Value = -1111.9
$ws.Range("N9").Value = -2737.5

$ws.Range("H18").Value = 567.5
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568

$ws.Range("H76").Value = 71434070
$ws.Range("I76").Value = 4376.3335
$ws.Range("J76").Value = 125006350
$ws.Range("K76").Value = 4376.3335
$ws.Range("L76").Value = 125006350
$ws.Range("M76").Value = -4061.3335
$ws.Range("N76").Value = -125006980

$ws.Range("H79").Value = 71434070
$ws.Range("I79").Value = 4376.3335
$ws.Range("J79").Value = 125006350
$ws.Range("K79").Value = 4376.3335
$ws.Range("L79").Value = 125006350
$ws.Range("M79").Value = -3284.3335
$ws.Range("N79").Value = -125008534

$ws.Range("H88").Value = 1754.75
$ws.Range("I88").Value = 988
$ws.Range("J88").Value = 1908.1
$ws.Range("K88").Value = 988
$ws.Range("L88").Value = 1908.1
$ws.Range("M88").Value = -582
$ws.Range("N88").Value = -2720.1

$ws.Range("H91").Value = 1754.75
$ws.Range("I91").Value = 988
$ws.Range("J91").Value = 1908.1
$ws.Range("K91").Value = 988
$ws.Range("L91").Value = 1908.1
$ws.Range("M91").Value = 416
$ws.Range("N91").Value = -4716.1

$ws.Range("H96").Value = 775.1539
$ws.Range("I96").Value = 701
$ws.Range("K96").Value = 2103
$ws.Range("M96").Value = -730

$ws.Range("H100").Value = 9214.789000000001
$ws.Range("I100").Value = 1732.5
$ws.Range("J100").Value = 12668.154
$ws.Range("K100").Value = 1732.5
$ws.Range("L100").Value = 12668.154
$ws.Range("M100").Value = -1191.5
$ws.Range("N100").Value = -13750.154

$ws.Range("H106").Value = 4947.5
$ws.Range("I106").Value = 4947.5
$ws.Range("K106").Value = 4947.5
$ws.Range("M106").Value = -4316.5

$ws.Range("H107").Value = 42247.082
$ws.Range("I107").Value = 44075.434
$ws.Range("K107").Value = 44075.434
$ws.Range("M107").Value = -42155.434

$ws.Range("H111").Value = 500497
$ws.Range("I111").Value = 999999
$ws.Range("J111").Value = 995
$ws.Range("K111").Value = 2999997
$ws.Range("L111").Value = 2985
$ws.Range("M111").Value = -2996930
$ws.Range("N111").Value = -9119

$ws.Range("H116").Value = 18871.5
$ws.Range("I116").Value = 7242
$ws.Range("K116").Value = 7242
$ws.Range("M116").Value = -3800

$ws.Range("H133").Value = 43581.08
$ws.Range("J133").Value = 43581.08
$ws.Range("L133").Value = 43581.08
$ws.Range("N133").Value = -53701.08

$ws.Range("H138").Value = 4366.22
$ws.Range("I138").Value = 2448.7144
$ws.Range("J138").Value = 4962.778
$ws.Range("K138").Value = 7346.1432
$ws.Range("L138").Value = 14888.334
$ws.Range("M138").Value = -2206.1432
$ws.Range("N138").Value = -25168.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7076.7334
$ws.Range("I2").Value = 9421.546
$ws.Range("K2").Value = 9421.546
$ws.Range("M2").Value = -9308.546

$ws.Range("H41").Value = 399
$ws.Range("I41").Value = 399
$ws.Range("K41").Value = 399
$ws.Range("M41").Value = 15

$ws.Range("H88").Value = 3343.3076
$ws.Range("I88").Value = 3378.4443
$ws.Range("K88").Value = 3378.4443
$ws.Range("M88").Value = -2972.4443

$ws.Range("H91").Value = 3343.3076
$ws.Range("I91").Value = 3378.4443
$ws.Range("K91").Value = 3378.4443
$ws.Range("M91").Value = -1974.4443

$ws.Range("H102").Value = 1600.8334
$ws.Range("I102").Value = 1673.6364
$ws.Range("K102").Value = 1673.6364
$ws.Range("M102").Value = -51.63640000000009

$ws.Range("H113").Value = 90398
$ws.Range("J113").Value = 90398
$ws.Range("L113").Value = 90398
$ws.Range("N113").Value = -99076

$ws.Range("H116").Value = 7076.7334
$ws.Range("I116").Value = 9421.546
$ws.Range("K116").Value = 9421.546
$ws.Range("M116").Value = -7127.546

$ws.Range("H132").Value = 4695.933
$ws.Range("I132").Value = 4695.933
$ws.Range("K132").Value = 14087.799
$ws.Range("M132").Value = -11557.799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7076.7334
$ws.Range("I3").Value = 9421.546
$ws.Range("K3").Value = 9421.546
$ws.Range("M3").Value = -9307.546

$ws.Range("H20").Value = 3978.7334
$ws.Range("I20").Value = 3608.889
$ws.Range("J20").Value = 4533.5
$ws.Range("K20").Value = 3608.889
$ws.Range("L20").Value = 4533.5
$ws.Range("M20").Value = -3361.889
$ws.Range("N20").Value = -5027.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 64161.145
$ws.Range("J20").Value = 64161.145
$ws.Range("L20").Value = 64161.145
$ws.Range("N20").Value = -64633.145

$ws.Range("H30").Value = 64161.145
$ws.Range("J30").Value = 64161.145
$ws.Range("L30").Value = 64161.145
$ws.Range("N30").Value = -64343.145

$ws.Range("H128").Value = 64161.145
$ws.Range("J128").Value = 64161.145
$ws.Range("L128").Value = 64161.145
$ws.Range("N128").Value = -74121.14499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 349.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 349.5
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -1272.5

$ws.Range("H44").Value = 300
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H131").Value = 2999.1614
$ws.Range("I131").Value = 1721.8572
$ws.Range("J131").Value = 3161.7273
$ws.Range("K131").Value = 5165.571599999999
$ws.Range("L131").Value = 9485.1819
$ws.Range("M131").Value = -125.5715999999993
$ws.Range("N131").Value = -19565.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 125014880
$ws.Range("I70").Value = 14401.2
$ws.Range("J70").Value = 333349000
$ws.Range("K70").Value = 14401.2
$ws.Range("L70").Value = 333349000
$ws.Range("M70").Value = -14131.2
$ws.Range("N70").Value = -333349540

$ws.Range("H73").Value = 125014880
$ws.Range("I73").Value = 14401.2
$ws.Range("J73").Value = 333349000
$ws.Range("K73").Value = 14401.2
$ws.Range("L73").Value = 333349000
$ws.Range("M73").Value = -13465.2
$ws.Range("N73").Value = -333350872

$ws.Range("H80").Value = 1821673.1
$ws.Range("J80").Value = 2501550.8
$ws.Range("L80").Value = 2501550.8
$ws.Range("N80").Value = -2503546.8

$ws.Range("H83").Value = 1821673.1
$ws.Range("J83").Value = 2501550.8
$ws.Range("L83").Value = 12507754
$ws.Range("N83").Value = -12517738

$ws.Range("H113").Value = 635133.4
$ws.Range("I113").Value = 1253395
$ws.Range("K113").Value = 1253395
$ws.Range("M113").Value = -1251225

$ws.Range("H132").Value = 44848.44
$ws.Range("I132").Value = 4110.8
$ws.Range("K132").Value = 12332.4
$ws.Range("M132").Value = -9802.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3024
$ws.Range("I68").Value = 2747
$ws.Range("K68").Value = 2747
$ws.Range("M68").Value = -1998

$ws.Range("H71").Value = 3024
$ws.Range("I71").Value = 2747
$ws.Range("K71").Value = 13735
$ws.Range("M71").Value = -9991

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1888.375
$ws.Range("I81").Value = 1701.9166
$ws.Range("K81").Value = 3403.8332
$ws.Range("M81").Value = -2342.8332

$ws.Range("H84").Value = 1888.375
$ws.Range("I84").Value = 1701.9166
$ws.Range("K84").Value = 17019.166
$ws.Range("M84").Value = -11715.166

$ws.Range("H132").Value = 43892.08
$ws.Range("I132").Value = 3090.1
$ws.Range("K132").Value = 9270.299999999999
$ws.Range("M132").Value = -6740.299999999999

$ws.Range("H136").Value = 9100109
$ws.Range("I136").Value = 11459718
$ws.Range("K136").Value = 34379154
$ws.Range("M136").Value = -34376604

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()
